$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 15 so the existing OUTA row shifts down to row 16;
# Excel copies the formatting of the row above into the newly inserted row.
$ws.Rows.Item(15).Insert()

# Row 15 becomes the new LDA entry
$ws.Cells.Item(15, 1).Value = "LDA"
$ws.Cells.Item(15, 2).Value = 1101
$ws.Cells.Item(15, 3).Value = "Cargar el registro desde memoria"

# Row 16 keeps the OUTA entry, but with updated binary code
$ws.Cells.Item(16, 2).Value = 1110

$ws.Range("B8").Select()
